$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 22251
$ws.Range("I9").Value = 28582.715
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 28582.715
$ws.Range("L9").Value = 90
$ws.Range("M9").Value = -28413.715
$ws.Range("N9").Value = -428
$ws.Range("H21").Value = 26208.143
$ws.Range("I21").Value = 29491
$ws.Range("K21").Value = 29491
$ws.Range("M21").Value = -29023
$ws.Range("H23").Value = 26208.143
$ws.Range("I23").Value = 29491
$ws.Range("K23").Value = 29491
$ws.Range("M23").Value = -29257
$ws.Range("H29").Value = 308.83334
$ws.Range("I29").Value = 270.6
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 811.8000000000001
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = -530.8000000000001
$ws.Range("N29").Value = -2062
$ws.Range("H32").Value = 973.2143
$ws.Range("I32").Value = 1226
$ws.Range("J32").Value = 832.7778
$ws.Range("K32").Value = 1226
$ws.Range("L32").Value = 832.7778
$ws.Range("M32").Value = -900
$ws.Range("N32").Value = -1484.7778
$ws.Range("H38").Value = 173.9
$ws.Range("I38").Value = 29.875
$ws.Range("J38").Value = 750
$ws.Range("K38").Value = 89.625
$ws.Range("L38").Value = 2250
$ws.Range("M38").Value = 282.375
$ws.Range("N38").Value = -2994
$ws.Range("H41").Value = 429.05554
$ws.Range("I41").Value = 608.3
$ws.Range("J41").Value = 205
$ws.Range("K41").Value = 608.3
$ws.Range("L41").Value = 205
$ws.Range("M41").Value = -168.3
$ws.Range("N41").Value = -1085
$ws.Range("H58").Value = 310
$ws.Range("I58").Value = 310
$ws.Range("K58").Value = 930
$ws.Range("M58").Value = -780
$ws.Range("H135").Value = 1086.4
$ws.Range("I135").Value = 422.10526
$ws.Range("J135").Value = 2233.818
$ws.Range("K135").Value = 3798.94734
$ws.Range("L135").Value = 20104.362
$ws.Range("M135").Value = -1263.94734
$ws.Range("N135").Value = -25174.362
$ws.Range("H137").Value = 4239.1
$ws.Range("I137").Value = 1240.1428
$ws.Range("J137").Value = 5853.923
$ws.Range("K137").Value = 3720.4284
$ws.Range("L137").Value = 17561.769
$ws.Range("M137").Value = -1170.4284
$ws.Range("N137").Value = -22661.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8507.409
$ws.Range("I45").Value = 9376.263000000001
$ws.Range("J45").Value = 3004.6667
$ws.Range("K45").Value = 9376.263000000001
$ws.Range("L45").Value = 3004.6667
$ws.Range("M45").Value = -8999.263000000001
$ws.Range("N45").Value = -3758.6667
$ws.Range("H122").Value = 2338
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 2014
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 6042
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -10942

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227
$ws.Range("H64").Value = 378.82352
$ws.Range("I64").Value = 220.4
$ws.Range("J64").Value = 444.83334
$ws.Range("K64").Value = 220.4
$ws.Range("L64").Value = 444.83334
$ws.Range("M64").Value = 4.599999999999994
$ws.Range("N64").Value = -894.83334
$ws.Range("H67").Value = 378.82352
$ws.Range("I67").Value = 220.4
$ws.Range("J67").Value = 444.83334
$ws.Range("K67").Value = 220.4
$ws.Range("L67").Value = 444.83334
$ws.Range("M67").Value = 559.6
$ws.Range("N67").Value = -2004.83334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 16500
$ws.Range("J28").Value = 16500
$ws.Range("L28").Value = 16500
$ws.Range("N28").Value = -16990
$ws.Range("H132").Value = 37043784
$ws.Range("I132").Value = 55563732
$ws.Range("J132").Value = 3895.3333
$ws.Range("K132").Value = 166691196
$ws.Range("L132").Value = 11685.9999
$ws.Range("M132").Value = -166688666
$ws.Range("N132").Value = -16745.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 623.3889
$ws.Range("I14").Value = 623.3889
$ws.Range("K14").Value = 1870.1667
$ws.Range("M14").Value = -1697.1667
$ws.Range("H23").Value = 47.94737
$ws.Range("I23").Value = 47.6
$ws.Range("J23").Value = 48.333332
$ws.Range("K23").Value = 142.8
$ws.Range("L23").Value = 144.999996
$ws.Range("M23").Value = 92.19999999999999
$ws.Range("N23").Value = -614.999996
$ws.Range("H107").Value = 580.8
$ws.Range("I107").Value = 450
$ws.Range("J107").Value = 668
$ws.Range("K107").Value = 1350
$ws.Range("L107").Value = 2004
$ws.Range("M107").Value = 570
$ws.Range("N107").Value = -5844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1524.05
$ws.Range("I122").Value = 1475.9166
$ws.Range("K122").Value = 4427.7498
$ws.Range("M122").Value = -1977.7498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1875
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -698
$ws.Range("N61").Value = -2604
$ws.Range("H97").Value = 12048
$ws.Range("J97").Value = 12048
$ws.Range("L97").Value = 12048
$ws.Range("N97").Value = -14030
$ws.Range("H113").Value = 1875
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -6540
$ws.Range("H122").Value = 253249.5
$ws.Range("I122").Value = 501499
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 1504497
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1502047
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 15505.182
$ws.Range("I132").Value = 26540
$ws.Range("J132").Value = 6309.5
$ws.Range("K132").Value = 79620
$ws.Range("L132").Value = 18928.5
$ws.Range("M132").Value = -77090
$ws.Range("N132").Value = -23988.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 283
$ws.Range("I107").Value = 266
$ws.Range("J107").Value = 295.75
$ws.Range("K107").Value = 798
$ws.Range("L107").Value = 887.25
$ws.Range("M107").Value = 1122
$ws.Range("N107").Value = -4727.25
$ws.Range("H122").Value = 11766209
$ws.Range("I122").Value = 18183162
$ws.Range("K122").Value = 54549486
$ws.Range("M122").Value = -54547036
$ws.Range("H136").Value = 1578.8928
$ws.Range("I136").Value = 789.26666
$ws.Range("J136").Value = 2490
$ws.Range("K136").Value = 2367.79998
$ws.Range("L136").Value = 7470
$ws.Range("M136").Value = 182.2000200000002
$ws.Range("N136").Value = -12570
